# feat: add 2022-Q1 data
#
# Before: 2 sheets  -> "2021-Q4", "总计"
# After:  3 sheets  -> "2021-Q4", "2022-Q1", "总计"
#   * "2022-Q1" is inserted between "2021-Q4" and "总计" and mirrors the
#     layout/headers of "2021-Q4" but with its own holding data.
#   * "总计" (the summary sheet) gains a new row for "2022-Q1" above the
#     existing "2021-Q4" row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right after "2021-Q4".
# ---------------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Item("2021-Q4")
$wsQ1 = $wb.Worksheets.Add($null, $wsQ4)
$wsQ1.Name = "2022-Q1"

# Sheet handles captured before the insert above can go stale (worksheet
# collection indices shift), so re-resolve everything we still need by
# name now that the sheet count/order is final.
$wsQ4 = $wb.Worksheets.Item("2021-Q4")
$wsQ1 = $wb.Worksheets.Item("2022-Q1")
$wsTotal = $wb.Worksheets.Item("总计")

# Copy header row + row 2 (values & styles) from "2021-Q4" as a starting
# point -- keeps the same column layout/formatting (bold/centered/bordered
# header, style on A2, etc.) without re-creating styles by hand.
$wsQ4.Range("B1:H1").Copy($wsQ1.Range("B1:H1"))
$wsQ4.Range("A2:H2").Copy($wsQ1.Range("A2:H2"))

# Now overwrite row 2 with the 2022-Q1 fund data. B2 (fund code) and
# D2:G2 hold numeric-looking values that are stored as text (like the
# source data), so mark them as Text before writing to keep the leading
# zero on the fund code and avoid the other columns being read as numbers.
$wsQ1.Range("B2").NumberFormat = "@"
$wsQ1.Range("D2:G2").NumberFormat = "@"
$wsQ1.Range("B2").Value = "006282"
$wsQ1.Range("C2").Value = "上投摩根欧洲动力策略股票（QDII）"
$wsQ1.Range("D2").Value = "0.48"
$wsQ1.Range("E2").Value = "89.68"
$wsQ1.Range("F2").Value = "1.87"
$wsQ1.Range("G2").Value = "0.0090"
$wsQ1.Range("H2").Value = 10

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: push the existing data row down to row 3
#    and add the new 2022-Q1 row in its place (row 2).
# ---------------------------------------------------------------------
$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3:D3"))
$wsTotal.Range("A3").Value = 1

$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.01

# Keep the originally active sheet selected (adding/editing sheets above
# shifts the active tab), matching the workbook's original view state.
$wsQ4.Activate()

